$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (column D) values
$ws.Range("D2").Value = "65.996.26"
$ws.Range("D3").Value = "2.957.64"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "572.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "164.61"
$ws.Range("D6").Style = "Normal"
$ws.Range("D9").Value = "2.956.33"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.68"
$ws.Range("D10").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000245"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.47"
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").Value = "65.890.71"
$ws.Range("D17").Value = "3.450.19"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.03"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "2.957.27"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "448.74"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.698"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.31"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.27"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.25"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.36"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.07"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.26"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.47"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.60"
$ws.Range("D31").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.118"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "27.40"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.974"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.74"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "47.88"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "49.16"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.99"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.305"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.83"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.49"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "384.29"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0352"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Value = "2.683.37"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "133.67"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.02"
$ws.Range("D50").Style = "Normal"

# Update Volume(1h) (column E) values
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("E3").Value = "  -1.91%  "
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("E5").Value = "  -2.11%  "
$ws.Range("E6").Value = "  +1.91%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  -0.36%  "
$ws.Range("E9").Value = "  -1.82%  "
$ws.Range("E10").Value = "  -2.34%  "
$ws.Range("E11").Value = "  -2.89%  "
$ws.Range("E12").Value = "  +1.38%  "
$ws.Range("E13").Value = "  -3.01%  "
$ws.Range("E14").Value = "  +2.33%  "
$ws.Range("E15").Value = "  -0.45%  "
$ws.Range("E16").Value = "  -0.05%  "
$ws.Range("E17").Value = "  -1.85%  "
$ws.Range("E18").Value = "  +2.28%  "
$ws.Range("E19").Value = "  +14.79%  "
$ws.Range("E20").Value = "  -1.91%  "
$ws.Range("E21").Value = "  -1.81%  "
$ws.Range("E22").Value = "  +1.42%  "
$ws.Range("E23").Value = "  -1.18%  "
$ws.Range("E24").Value = "  -0.13%  "
$ws.Range("E25").Value = "  -0.98%  "
$ws.Range("E26").Value = "  -0.47%  "
$ws.Range("E27").Value = "  -5.63%  "
$ws.Range("E28").Value = "  +0.16%  "
$ws.Range("E29").Value = "  +2.78%  "
$ws.Range("E30").Value = "  +4.77%  "
$ws.Range("E31").Value = "  -0.32%  "
$ws.Range("E32").Value = "  -4.90%  "
$ws.Range("E33").Value = "  +6.38%  "
$ws.Range("E34").Value = "  +0.84%  "
$ws.Range("E35").Value = "  -0.16%  "
$ws.Range("E36").Value = "  -2.04%  "
$ws.Range("E37").Value = "  -1.71%  "
$ws.Range("E38").Value = "  +10.23%  "
$ws.Range("E39").Value = "  -1.50%  "
$ws.Range("E40").Value = "  -8.25%  "
$ws.Range("E41").Value = "  -1.78%  "
$ws.Range("E42").Value = "  -1.16%  "
$ws.Range("E43").Value = "  -4.78%  "
$ws.Range("E44").Value = "  +0.52%  "
$ws.Range("E45").Value = "  -1.61%  "
$ws.Range("E46").Value = "  -1.00%  "
$ws.Range("E47").Value = "  -4.11%  "
$ws.Range("E48").Value = "  -0.52%  "
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("E50").Value = "  +1.11%  "
$ws.Range("E51").Value = "  +1.83%  "
